# "Classes específicas para cidade"
#
# The source edit mainly re-flows existing sentences into additional
# <w:r> runs and sprinkles <w:proofErr> spell/grammar markers around
# loan-words (GitHub, Xamarin, Android, RelativeLayout, canvas, hash,
# foreach, IEnumerable, GetEnumerator, StreamReader, Assets, Canvas,
# openFileInput/openFileOutput, etc). None of that changes the visible
# text or its formatting - it is cosmetic, Word-internal bookkeeping
# that is produced by Word's proofing engine, not something the Word
# object model exposes a way to author directly.
#
# The one substantive, visible change in the diff is in the bullet
# "28/11 às 11h00min: classes Aresta<T> e PesoCidades." - the stray
# "<T>" generic-type marker is dropped from "Aresta<T>" (leaving just
# "Aresta") and the resulting double space collapses to one, giving:
#   "28/11 às 11h00min: classes Aresta e PesoCidades."
# The "Aresta" / "PesoCidades" words keep their italics and "e" keeps
# the surrounding roman (non-italic) text, matching the target markup.
#
# We scope every Find to the specific sentence (via a Range captured
# from an unambiguous, one-off match) so we never touch any of the
# several other " e " occurrences elsewhere in the document.

$d = $word.ActiveDocument

$scope = $d.Content
$scope.Find.Text = "28/11 às 11h00min: classes Aresta<T>  e PesoCidades."
$null = $scope.Find.Execute()
$sentenceStart = $scope.Start
$sentenceEnd = $scope.End

# "Aresta<T> " (the italic run) -> "Aresta " (drop the <T> generic marker)
$r1 = $d.Range($sentenceStart, $sentenceEnd)
$r1.Find.Execute("Aresta<T> ", $true, $false, $false, $false, $false,
                  $true, 1, $false, "Aresta ", 2)

# " e " (the roman run) -> "e " (the leading space was already supplied
# by the trailing space kept in the italic "Aresta " run above)
$r2 = $d.Range($sentenceStart, $sentenceEnd)
$r2.Find.Execute(" e ", $true, $false, $false, $false, $false,
                  $true, 1, $false, "e ", 2)

# Word keeps the "_GoBack" last-edit bookmark pinned at the point of the
# most recent edit; after collapsing "Aresta<T>  " -> "Aresta " that
# point sits right after the word "Aresta " (before "e PesoCidades."),
# matching the new bookmark position shown in the diff.
$r3 = $d.Range($sentenceStart, $sentenceEnd)
$r3.Find.Text = "Aresta "
$null = $r3.Find.Execute()
$editPoint = $d.Range($r3.End, $r3.End)
$d.Bookmarks.Add("_GoBack", $editPoint)
